$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.706
$ws.Range("A8").Value = -22.305
$ws.Range("A10").Value = -21.889
$ws.Range("D11").Value = -7.264
$ws.Range("A12").Value = -21.474
$ws.Range("D12").Value = -6.873
$ws.Range("D15").Value = -8.317
$ws.Range("D17").Value = -8.193000000000001
$ws.Range("A18").Value = -22.095
$ws.Range("A25").Value = -21.654
$ws.Range("D26").Value = -7.255000000000001
$ws.Range("D27").Value = -7.728
$ws.Range("D28").Value = -8.047999999999998
$ws.Range("D32").Value = -7.302000000000001
$ws.Range("A37").Value = -20.317
$ws.Range("D37").Value = -8.108000000000001
$ws.Range("D41").Value = -8.061
$ws.Range("D47").Value = -7.452
$ws.Range("D51").Value = -8.335000000000001
$ws.Range("A55").Value = -22.204
$ws.Range("D65").Value = -7.741
$ws.Range("A68").Value = -21.565
$ws.Range("D73").Value = -8.018000000000001
$ws.Range("A77").Value = -20.236
$ws.Range("A78").Value = -19.836
$ws.Range("A79").Value = -21.205
$ws.Range("A80").Value = -20.17
$ws.Range("A81").Value = -21.697
$ws.Range("A82").Value = -22.209
$ws.Range("A84").Value = -21.849
$ws.Range("D84").Value = -8.161000000000001
$ws.Range("D85").Value = -8.790000000000001
$ws.Range("D89").Value = -8.292
$ws.Range("D93").Value = -6.872999999999999
$ws.Range("D95").Value = -7.557
$ws.Range("D98").Value = -7.204000000000001
$ws.Range("D99").Value = -8.272000000000002
$ws.Range("A101").Value = -20.478
$ws.Range("D101").Value = -7.834999999999999
$ws.Range("A102").Value = -20.431
$ws.Range("D102").Value = -8.108000000000001
